$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$sh = $s.Shapes.Item("TextBox 22")
$tr = $sh.TextFrame.TextRange

# Paragraph 5 is "Fit & Predict" (bulleted). Insert two new bulleted
# paragraphs ("Scaling", "Encoding") immediately before it, reusing its
# paragraph formatting (bullet char/font, indent) for the new paragraphs.
$fitPredict = $tr.Paragraphs(5, 1)
$fitPredict.InsertBefore("Scaling" + [char]13 + "Encoding" + [char]13) | Out-Null

# After the insert, "Fit & Predict" has shifted down to paragraph 7.
# Insert a new bulleted paragraph ("PCA") right after it.
$fitPredict = $tr.Paragraphs(7, 1)
$fitPredict.InsertAfter([char]13 + "PCA") | Out-Null

# The textbox has spAutoFit; grow it to fit the three extra lines
# (2031325 EMU -> 2862322 EMU, i.e. 159.9469pt -> 225.3797pt).
$sh.Height = 225.3797
